$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Аксессуары"
$ws.Range("A3").Value = "Декор"
$ws.Range("A4").Value = "Книги"
$ws.Range("A5").Value = "Косметика"
$ws.Range("A6").Value = "Кулинария"
$ws.Range("A7").Value = "Игры"
$ws.Range("A8").Value = "Одежда"
$ws.Range("A9").Value = "Спорт"
$ws.Range("A10").Value = "Хобби"
$ws.Range("A11").Value = "Гаджеты"

$ws.Range("A11").Select()
